$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "model_4_0_0"
$ws.Range("B2").Value = -0.4469409842161907
$ws.Range("C2").Value = -0.04918934452759949
$ws.Range("D2").Value = -5.404944135831724
$ws.Range("E2").Value = -0.6835525857389761
$ws.Range("F2").Value = 1.601337432861328
$ws.Range("G2").Value = 2.265443325042725
$ws.Range("H2").Value = 2.090804100036621
$ws.Range("I2").Value = 2.183259963989258

$ws.Range("A3").Value = "model_4_0_1"
$ws.Range("B3").Value = -0.245215902882957
$ws.Range("C3").Value = 0.07559084591087306
$ws.Range("D3").Value = -4.817674746363592
$ws.Range("E3").Value = -0.5039948414662416
$ws.Range("F3").Value = 1.378087162971497
$ws.Range("G3").Value = 1.99601411819458
$ws.Range("H3").Value = 1.89909827709198
$ws.Range("I3").Value = 1.950406432151794

$ws.Range("A4").Value = "model_4_0_2"
$ws.Range("B4").Value = 0.1491667311575263
$ws.Range("C4").Value = 0.3762570554252721
$ws.Range("D4").Value = -3.680254514736075
$ws.Range("E4").Value = -0.1042273170938899
$ws.Range("F4").Value = 0.9416216611862183
$ws.Range("G4").Value = 1.346805810928345
$ws.Range("H4").Value = 1.527803421020508
$ws.Range("I4").Value = 1.431981205940247

$ws.Range("A5").Value = "model_4_0_21"
$ws.Range("B5").Value = 0.4752361594400787
$ws.Range("C5").Value = 0.3004030442834689
$ws.Range("D5").Value = -4.227368003641288
$ws.Range("E5").Value = -0.2359013007504949
$ws.Range("F5").Value = 0.5807589292526245
$ws.Range("G5").Value = 1.510592341423035
$ws.Range("H5").Value = 1.706400990486145
$ws.Range("I5").Value = 1.602738380432129

$ws.Range("A6").Value = "model_4_0_20"
$ws.Range("B6").Value = 0.4753079634901229
$ws.Range("C6").Value = 0.3007357692109894
$ws.Range("D6").Value = -4.225066096611656
$ws.Range("E6").Value = -0.2353345622279057
$ws.Range("F6").Value = 0.5806794762611389
$ws.Range("G6").Value = 1.509873747825623
$ws.Range("H6").Value = 1.705649495124817
$ws.Range("I6").Value = 1.602003335952759

$ws.Range("A7").Value = "model_4_0_19"
$ws.Range("B7").Value = 0.4756795257024039
$ws.Range("C7").Value = 0.3033036057699519
$ws.Range("D7").Value = -4.234113784528622
$ws.Range("E7").Value = -0.2341426202555477
$ws.Range("F7").Value = 0.5802683234214783
$ws.Range("G7").Value = 1.504329323768616
$ws.Range("H7").Value = 1.708603024482727
$ws.Range("I7").Value = 1.600457549095154

$ws.Range("A8").Value = "model_4_0_18"
$ws.Range("B8").Value = 0.4773026590226335
$ws.Range("C8").Value = 0.3059693294911021
$ws.Range("D8").Value = -4.214797434612629
$ws.Range("E8").Value = -0.2295050169766086
$ws.Range("F8").Value = 0.5784719586372375
$ws.Range("G8").Value = 1.498573184013367
$ws.Range("H8").Value = 1.702297449111938
$ws.Range("I8").Value = 1.594443321228027

$ws.Range("A9").Value = "model_4_0_17"
$ws.Range("B9").Value = 0.4782241404492402
$ws.Range("C9").Value = 0.3084251855817413
$ws.Range("D9").Value = -4.208447730777017
$ws.Range("E9").Value = -0.2265879156132748
$ws.Range("F9").Value = 0.5774521231651306
$ws.Range("G9").Value = 1.493270516395569
$ws.Range("H9").Value = 1.700224637985229
$ws.Range("I9").Value = 1.590660572052002

$ws.Range("A10").Value = "model_4_0_16"
$ws.Range("B10").Value = 0.4823600247453342
$ws.Range("C10").Value = 0.3089012496183828
$ws.Range("D10").Value = -4.096173589896254
$ws.Range("E10").Value = -0.212868970808445
$ws.Range("F10").Value = 0.572874903678894
$ws.Range("G10").Value = 1.492242693901062
$ws.Range("H10").Value = 1.66357433795929
$ws.Range("I10").Value = 1.572869539260864

$ws.Range("A11").Value = "model_4_0_15"
$ws.Range("B11").Value = 0.4850507552359519
$ws.Range("C11").Value = 0.3150867357077864
$ws.Range("D11").Value = -4.072498595822275
$ws.Range("E11").Value = -0.2046116551664421
$ws.Range("F11").Value = 0.5698970556259155
$ws.Range("G11").Value = 1.478886723518372
$ws.Range("H11").Value = 1.655845999717712
$ws.Range("I11").Value = 1.562161326408386

$ws.Range("A12").Value = "model_4_0_14"
$ws.Range("B12").Value = 0.4856520956371132
$ws.Range("C12").Value = 0.3173974103583042
$ws.Range("D12").Value = -4.073668996871977
$ws.Range("E12").Value = -0.2027133783826636
$ws.Range("F12").Value = 0.5692315697669983
$ws.Range("G12").Value = 1.473897337913513
$ws.Range("H12").Value = 1.656228065490723
$ws.Range("I12").Value = 1.559699654579163

$ws.Range("A13").Value = "model_4_0_13"
$ws.Range("B13").Value = 0.4866020189877961
$ws.Range("C13").Value = 0.3213710370959386
$ws.Range("D13").Value = -4.07831827316657
$ws.Range("E13").Value = -0.1997613887951444
$ws.Range("F13").Value = 0.56818026304245
$ws.Range("G13").Value = 1.465317368507385
$ws.Range("H13").Value = 1.657745838165283
$ws.Range("I13").Value = 1.555871367454529

$ws.Range("A14").Value = "model_4_0_12"
$ws.Range("B14").Value = 0.4883117076712583
$ws.Range("C14").Value = 0.3300835832279264
$ws.Range("D14").Value = -4.096381049882757
$ws.Range("E14").Value = -0.1942215934177445
$ws.Range("F14").Value = 0.5662881731987
$ws.Range("G14").Value = 1.446505069732666
$ws.Range("H14").Value = 1.663642048835754
$ws.Range("I14").Value = 1.548687219619751

$ws.Range("A15").Value = "model_4_0_11"
$ws.Range("B15").Value = 0.4965965316860412
$ws.Range("C15").Value = 0.3558906665504805
$ws.Range("D15").Value = -4.060673557935568
$ws.Range("E15").Value = -0.167243232152561
$ws.Range("F15").Value = 0.5571193099021912
$ws.Range("G15").Value = 1.39078164100647
$ws.Range("H15").Value = 1.651985883712769
$ws.Range("I15").Value = 1.513701438903809

$ws.Range("A16").Value = "model_4_0_10"
$ws.Range("B16").Value = 0.4992950925907543
$ws.Range("C16").Value = 0.3410794655955914
$ws.Range("D16").Value = -3.843486047120784
$ws.Range("E16").Value = -0.1545715891703638
$ws.Range("F16").Value = 0.5541328191757202
$ws.Range("G16").Value = 1.422762393951416
$ws.Range("H16").Value = 1.581087946891785
$ws.Range("I16").Value = 1.497268557548523

$ws.Range("A17").Value = "model_4_0_9"
$ws.Range("B17").Value = 0.502321120644556
$ws.Range("C17").Value = 0.3615855782561286
$ws.Range("D17").Value = -3.756325756400034
$ws.Range("E17").Value = -0.1261707418013347
$ws.Range("F17").Value = 0.5507838726043701
$ws.Range("G17").Value = 1.378484845161438
$ws.Range("H17").Value = 1.552635788917542
$ws.Range("I17").Value = 1.460437774658203

$ws.Range("A18").Value = "model_4_0_4"
$ws.Range("B18").Value = 0.5407184309919975
$ws.Range("C18").Value = 0.7318735216070309
$ws.Range("D18").Value = -3.428040368763111
$ws.Range("E18").Value = 0.2391189596396813
$ws.Range("F18").Value = 0.5082893371582031
$ws.Range("G18").Value = 0.5789473056793213
$ws.Range("H18").Value = 1.44547176361084
$ws.Range("I18").Value = 0.9867237210273743

$ws.Range("A19").Value = "model_4_0_3"
$ws.Range("B19").Value = 0.5468175804075355
$ws.Range("C19").Value = 0.7918706708724604
$ws.Range("D19").Value = -3.125347516395973
$ws.Range("E19").Value = 0.3278616268183112
$ws.Range("F19").Value = 0.501539409160614
$ws.Range("G19").Value = 0.4493995308876038
$ws.Range("H19").Value = 1.346661806106567
$ws.Range("I19").Value = 0.871640682220459

$ws.Range("A20").Value = "model_4_0_5"
$ws.Range("B20").Value = 0.5680471648867598
$ws.Range("C20").Value = 0.6780351653490481
$ws.Range("D20").Value = -3.631805924194773
$ws.Range("E20").Value = 0.1675242724788083
$ws.Range("F20").Value = 0.4780445098876953
$ws.Range("G20").Value = 0.6951968669891357
$ws.Range("H20").Value = 1.511988043785095
$ws.Range("I20").Value = 1.079568862915039

$ws.Range("A21").Value = "model_4_0_7"
$ws.Range("B21").Value = 0.5814784692165009
$ws.Range("C21").Value = 0.5785695699171833
$ws.Range("D21").Value = -3.430415802120319
$ws.Range("E21").Value = 0.1037029044906984
$ws.Range("F21").Value = 0.4631800055503845
$ws.Range("G21").Value = 0.9099661707878113
$ws.Range("H21").Value = 1.446247100830078
$ws.Range("I21").Value = 1.162333607673645

$ws.Range("A22").Value = "model_4_0_8"
$ws.Range("B22").Value = 0.5824534393935368
$ws.Range("C22").Value = 0.5628807347570158
$ws.Range("D22").Value = -3.324092744304412
$ws.Range("E22").Value = 0.1024682123114927
$ws.Range("F22").Value = 0.462101012468338
$ws.Range("G22").Value = 0.9438420534133911
$ws.Range("H22").Value = 1.411539435386658
$ws.Range("I22").Value = 1.163934707641602

$ws.Range("A23").Value = "model_4_0_6"
$ws.Range("B23").Value = 0.5887930217554078
$ws.Range("C23").Value = 0.6079355755216829
$ws.Range("D23").Value = -3.205099266256981
$ws.Range("E23").Value = 0.1562788591259082
$ws.Range("F23").Value = 0.4550849497318268
$ws.Range("G23").Value = 0.8465580940246582
$ws.Range("H23").Value = 1.372695565223694
$ws.Range("I23").Value = 1.094152212142944
